$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 333. This shifts the existing rows 333-413
# down to 334-414 and copies formatting (e.g. the date style on column D)
# from the row above, matching the target workbook.
$ws.Rows(333).Insert()

# Populate the newly inserted row 333 with its data.
$ws.Range("A333").Value = 5
$ws.Range("B333").Value = "Macroferia Regional de Talca"
$ws.Range("C333").Value = "Maule"
$ws.Range("D333").Value = 44855
$ws.Range("E333").Value = 7
$ws.Range("F333").Value = 100114014
$ws.Range("G333").Value = "Betarraga"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Segunda"
$ws.Range("J333").Value = 4000
$ws.Range("K333").Value = 800
$ws.Range("L333").Value = 800
$ws.Range("M333").Value = 800
$ws.Range("N333").Value = "`$/paquete 5 unidades"
$ws.Range("O333").Value = "Región del Maule"
$ws.Range("P333").Value = 160
$ws.Range("Q333").Value = 5
$ws.Range("R333").Value = "Hortaliza"
